$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-45 from 45190 to 45192
foreach ($row in 2..45) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
